$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Append a new row (row 33) of data, mirroring the existing rows above it.
$ws.Range("A33").Value = 30
$ws.Range("B33").Value = "5：13-6：09"
$ws.Range("C33").Value = "翻了翻C语言"
$ws.Range("E33").Value = "今天七月半，上午烧纸，下午挂水，摸鱼了（明天科二，鸽了）"

# Match the author's final selection on the newly added cell.
$ws.Range("E33").Select()
